$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New "Componentes Azure" comparison block (columns N:R, rows 15-22)
# ---------------------------------------------------------------------------
$ws.Range("N15").Value = "Componentes Azure"
$ws.Range("Q15").Value = "Núcleos"
$ws.Range("R15").Value = "Cantidad"

$ws.Range("N16").Value = "D16ds v5 16 vCPU"
$ws.Range("P16").Value = 87.6
$ws.Range("Q16").Value = 16
$ws.Range("R16").Value = 3

$ws.Range("N17").Value = "Disco 2 x 300 GiB SSD NVMe"
$ws.Range("P17").Value = 600
$ws.Range("R17").Value = 6

# N20/O20 mirror the "Nuevos tiempos de servicio" header (H20) formatting
$ws.Range("H20").Copy()
$ws.Range("N20:O20").PasteSpecial(-4122)
$ws.Range("O20").Value = "Nuevos tiempos de servicio"

$ws.Range("O21").Value = "Tpo cpu = (Nn x Ipb) /(Nb * Ipn cpu)*Tsb"
$ws.Range("P21").Formula = '=($C$6*$B$6)/(Q16*P16)*$B$2'

$ws.Range("O22").Value = "Tpo disco = Ipb x Tsb /Ipn"
$ws.Range("P22").Formula = '=$B$3*$B$7/P17'

# ---------------------------------------------------------------------------
# Existing formulas reworked to use absolute references / corrected formula
# ---------------------------------------------------------------------------
$ws.Range("I21").Formula = '=($C$6*$B$6)/(J16*I16)*$B$2'
$ws.Range("I22").Formula = '=$B$3*$B$7/I17'

# ---------------------------------------------------------------------------
# Row 25: replaced "Probabilidad disco" / literal 0.954545455 with a
# computed "tres" cost formula (duplicated for the on-prem vs Azure option)
# ---------------------------------------------------------------------------
$ws.Range("E25").Value = "tres "
$ws.Range("F25").Formula = '=($B$10*I21)+($B$11*I22)+($B$12*$B$4)'
$ws.Range("F25").Style = "Normal"
$ws.Range("L25").Value = "tres "
$ws.Range("M25").Formula = '=($B$10*P21)+($B$11*P22)+($B$12*$B$4)'
$ws.Range("M25").Style = "Normal"

# ---------------------------------------------------------------------------
# New pricing columns for the "Componentes elegidos" table (rows 28-37)
# ---------------------------------------------------------------------------
$ws.Range("I28").Value = "Precio"
$ws.Range("J28").Value = "Precio Total"

$ws.Range("I29").Value = 425
$ws.Range("J29").Formula = '=I29*E29'
$ws.Range("K29").Value = "€"
$ws.Range("K29").Font.Bold = $true

$ws.Range("I30").Value = 180
$ws.Range("J30").Formula = '=I30*E30'
$ws.Range("K30").Value = "€"
$ws.Range("K30").Font.Bold = $true

$ws.Range("I31").Value = 135
$ws.Range("J31").Formula = '=I31*E29'
$ws.Range("K31").Value = "€"
$ws.Range("K31").Font.Bold = $true

$ws.Range("I32").Value = 17
$ws.Range("J32").Formula = '=I32*E32'
$ws.Range("K32").Value = "€"
$ws.Range("K32").Font.Bold = $true

$ws.Range("I33").Value = 690
$ws.Range("J33").Formula = '=I33*E33'
$ws.Range("K33").Value = "€"
$ws.Range("K33").Font.Bold = $true

$ws.Range("I34").Value = 42
$ws.Range("J34").Formula = '=I34*E34'
$ws.Range("K34").Value = "€"
$ws.Range("K34").Font.Bold = $true

$ws.Range("I35").Value = 180
$ws.Range("J35").Formula = '=I35*E29'
$ws.Range("K35").Value = "€"
$ws.Range("K35").Font.Bold = $true

$ws.Range("J37").Formula = '=SUM(J29:J35)'
$ws.Range("K37").Value = "€"

# ---------------------------------------------------------------------------
# Conditional formatting on the grand total (green <=5400 / red >5400)
# ---------------------------------------------------------------------------
$fcGreen = $ws.Range("J37").FormatConditions.Add(2, 3, '$J$37<=5400')
$fcGreen.Interior.Color = 5296274
$fcRed = $ws.Range("J37").FormatConditions.Add(2, 3, '$J$37>5400')
$fcRed.Interior.Color = 255
$fcRed.SetFirstPriority()

# ---------------------------------------------------------------------------
# Misc view-state change
# ---------------------------------------------------------------------------
$ws.Range("I21").Select()
